$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached footer placeholder text (date + slide number) on the
#    slide master and on every slide layout. These are the auto-generated
#    <a:fld> caches PowerPoint re-writes whenever the deck is saved again.
# ---------------------------------------------------------------------------

$oldDate = "26-12-2017"
$newDate = "06/02/2018"

$oldSlideNum = [string]([char]0x2039) + "n" + [string][char]0xBA + [string][char]0x203A
$newSlideNum = [string]([char]0x2039) + "#" + [string][char]0x203A

function Update-FooterFields($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                $curText = $tr.Text
                if ($curText -eq $oldDate) {
                    $full = $tr.Characters(1, $tr.Length)
                    $full.Text = $newDate
                } elseif ($curText -eq $oldSlideNum) {
                    $full = $tr.Characters(1, $tr.Length)
                    $full.Text = $newSlideNum
                }
            }
        }
    }
}

$master = $p.Slides.Item(1).Master
Update-FooterFields $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-FooterFields $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Rename the "Pulse Shaper" blocks to "Electrical Filter" in the receiver
#    diagram on slide 1.
# ---------------------------------------------------------------------------

function Update-PulseShaper($shape) {
    if ($shape.Type -eq 6) {
        for ($j = 1; $j -le $shape.GroupItems.Count; $j++) {
            Update-PulseShaper $shape.GroupItems.Item($j)
        }
        return
    }
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "Pulse Shaper") {
                $full = $tr.Characters(1, $tr.Length)
                $full.Text = "Electrical Filter"
            }
        }
    }
}

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    Update-PulseShaper $slide.Shapes.Item($i)
}
